$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Position', 'First Name', 'Last Name', 'Shift', 'Location', 'Station')
    ,@('Bin Filler', ' Doris', 'Reynolds', 'Day', 'Presize', 'Winfield')
    ,@('Bin Filler', ' Zabada', 'Mohammed', 'Afternoon', 'Presize', 'Winfield')
    ,@('Bin Filler', ' Sukhwinder', 'Gakhal', 'Afternoon', 'Presize', 'Winfield')
    ,@('Bin Filler', ' Kuldip', 'Buttar', 'Day', 'Presize', 'Winfield')
    ,@('Forklift', ' George', 'Dunn', 'Day', 'Presize', 'Winfield')
    ,@('Forklift', ' Don', 'Coles', 'Day', 'Presize', 'Winfield')
    ,@('Forklift', ' George C', 'Brown', 'Day', 'Presize', 'Winfield')
    ,@('Forklift', ' Ron', 'Engene', 'Afternoon', 'Presize', 'Winfield')
    ,@('Forklift', ' Jerry', 'Engene', 'Afternoon', 'Presize', 'Winfield')
    ,@('Forklift', ' Fay', 'Lee', 'Afternoon', 'Presize', 'Winfield')
    ,@('Line Operator', ' Parveen', 'Gopal', 'Day', 'Presize', 'Winfield')
    ,@('Line Operator', ' Lori', 'Carter', 'Afternoon', 'Presize', 'Winfield')
    ,@('QC', ' Isabel', 'Roseen', 'Day', 'Presize', 'Winfield')
    ,@('QC', ' Wendy', 'Casorso', 'Afternoon', 'Presize', 'Winfield')
    ,@('Non Rotational', ' Elaine', 'Roseen', 'Day', 'Presize', 'Winfield')
    ,@('Non Rotational', ' Janeanne', 'Reiswig', 'Day', 'Presize', 'Winfield')
    ,@('Non Rotational', ' Sandra', 'Martin', 'Afternoon', 'Presize', 'Winfield')
    ,@('Non Rotational', ' Brian', 'High', 'Day', 'Presize', 'Winfield')
    ,@('Non Rotational', ' Joyce', 'Salga', 'Afternoon', 'Presize', 'Winfield')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
